$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.934.33"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "1.847.49"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  +0.43%  "
Set-TextValue $ws "D5" "310.01"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("E6").Value = "  +0.35%  "
Set-TextValue $ws "D7" "0.4709"
$ws.Range("E7").Value = "  +3.85%  "
Set-TextValue $ws "D8" "0.3666"
$ws.Range("E8").Value = "  +1.99%  "
Set-TextValue $ws "D9" "0.07164"
Set-TextValue $ws "D10" "0.9274"
$ws.Range("E10").Value = "  +3.76%  "
Set-TextValue $ws "D11" "19.61"
$ws.Range("E11").Value = "  +1.58%  "
Set-TextValue $ws "D12" "0.07712"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "1.804.42"
$ws.Range("E13").Value = "  -1.15%  "
Set-TextValue $ws "D14" "5.285"
$ws.Range("E14").Value = "  +0.48%  "
Set-TextValue $ws "D15" "6.412"
$ws.Range("E15").Value = "  +1.50%  "
Set-TextValue $ws "D16" "88.43"
$ws.Range("E16").Value = "  +3.69%  "
$ws.Range("E17").Value = "  +0.45%  "
Set-TextValue $ws "D18" "0.000008639"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "26.961.15"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").Value = "  +2.11%  "
Set-TextValue $ws "D22" "5.028"
$ws.Range("E22").Value = "  +1.34%  "
Set-TextValue $ws "D23" "10.62"
$ws.Range("E23").Value = "  +1.12%  "
Set-TextValue $ws "D24" "1.935"
$ws.Range("E24").Value = "  -1.36%  "
Set-TextValue $ws "D25" "152.05"
$ws.Range("E25").Value = "  +0.26%  "
Set-TextValue $ws "D26" "18.26"
$ws.Range("E26").Value = "  +2.52%  "
Set-TextValue $ws "D27" "2.021"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("E28").Value = "  +1.80%  "
Set-TextValue $ws "D29" "4.883"
$ws.Range("E29").Value = "  +0.90%  "
Set-TextValue $ws "D30" "0.08865"
$ws.Range("E30").Value = "  +1.73%  "
Set-TextValue $ws "D31" "3.218"
$ws.Range("E31").Value = "  +2.93%  "
Set-TextValue $ws "D32" "1.180"
$ws.Range("E32").Value = "  +6.03%  "
Set-TextValue $ws "D33" "0.7485"
$ws.Range("E33").Value = "  -0.43%  "
Set-TextValue $ws "D34" "2.781"
$ws.Range("E34").Value = "  +1.32%  "
Set-TextValue $ws "D35" "4.482"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("E36").Value = "  +1.30%  "
Set-TextValue $ws "D37" "0.01943"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D38" "0.05212"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D39" "2.960"
$ws.Range("E39").Value = "  +1.49%  "
Set-TextValue $ws "D40" "0.5231"
$ws.Range("E40").Value = "  +2.38%  "
Set-TextValue $ws "D41" "6.979"
$ws.Range("E41").Value = "  +3.02%  "
Set-TextValue $ws "D42" "0.1521"
$ws.Range("E42").Value = "  +0.64%  "
Set-TextValue $ws "D43" "8.172"
$ws.Range("E43").Value = "  +1.66%  "
Set-TextValue $ws "D44" "10.50"
$ws.Range("E44").Value = "  +5.21%  "
Set-TextValue $ws "D45" "0.4711"
$ws.Range("E45").Value = "  +0.12%  "
Set-TextValue $ws "D46" "1.007"
$ws.Range("E46").Value = "  +0.48%  "
Set-TextValue $ws "D47" "100.95"
$ws.Range("E47").Value = "  +0.87%  "
Set-TextValue $ws "D48" "1.599"
$ws.Range("E48").Value = "  +1.63%  "
Set-TextValue $ws "D49" "65.74"
$ws.Range("E49").Value = "  +2.82%  "
Set-TextValue $ws "D50" "0.06039"
$ws.Range("E50").Value = "  +0.94%  "
Set-TextValue $ws "D51" "0.8984"
$ws.Range("E51").Value = "  +6.06%  "
